$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated save_data: column G ("K" - strikeouts) replaces the old Strike# values
# with freshly calculated strikeout counts (std/mean recalculated, s_vals rewritten).
$newK = @{
    2  = 5
    3  = 6
    4  = 6
    5  = 6
    6  = 10
    7  = 2
    8  = 1
    9  = 4
    10 = 3
    11 = 4
    12 = 2
    13 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
